$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Resize columns L, M, N (12,13,14) to width 25 (ColumnWidth = 25 - 5/6 yields raw OOXML width 25)
$ws.Columns.Item(12).ColumnWidth = 25 - 5/6
$ws.Columns.Item(13).ColumnWidth = 25 - 5/6
$ws.Columns.Item(14).ColumnWidth = 25 - 5/6

# Row 2
$ws.Cells.Item(2,4).Value = "'2026-02-07"
$ws.Cells.Item(2,5).Value = "2026-02-07 00:17:58"
$ws.Cells.Item(2,6).Value = "https://www.meteo.cat/observacions/xema/dades?codi=YT&dia=2026-02-07T09:00Z"
$ws.Cells.Item(2,7).Value = "sense dades"
$ws.Cells.Item(2,8).Value = "sense dades"
$ws.Cells.Item(2,9).Value = "sense dades"
$ws.Cells.Item(2,10).ClearContents()
$ws.Cells.Item(2,11).Value = "sense dades"
$ws.Cells.Item(2,12).ClearContents()
$ws.Cells.Item(2,13).Value = "sense dades sense dades"
$ws.Cells.Item(2,14).Value = "sense dades sense dades"
$ws.Cells.Item(2,15).Value = "sense dades"

# Row 3
$ws.Cells.Item(3,4).Value = "'2026-02-07"
$ws.Cells.Item(3,5).Value = "2026-02-07 00:18:00"
$ws.Cells.Item(3,6).Value = "https://www.meteo.cat/observacions/xema/dades?codi=Z1&dia=2026-02-07T09:00Z"
$ws.Cells.Item(3,7).Value = "sense dades"
$ws.Cells.Item(3,8).Value = "sense dades"
$ws.Cells.Item(3,9).Value = "sense dades"
$ws.Cells.Item(3,10).ClearContents()
$ws.Cells.Item(3,11).Value = "sense dades"
$ws.Cells.Item(3,12).Value = "sense dades sense dades"
$ws.Cells.Item(3,13).Value = "sense dades sense dades"
$ws.Cells.Item(3,14).Value = "sense dades sense dades"
$ws.Cells.Item(3,15).Value = "sense dades"

# Row 4
$ws.Cells.Item(4,4).Value = "'2026-02-07"
$ws.Cells.Item(4,5).Value = "2026-02-07 00:18:02"
$ws.Cells.Item(4,6).Value = "https://www.meteo.cat/observacions/xema/dades?codi=DN&dia=2026-02-07T09:00Z"
$ws.Cells.Item(4,7).ClearContents()
$ws.Cells.Item(4,8).Value = "sense dades"
$ws.Cells.Item(4,9).Value = "sense dades"
$ws.Cells.Item(4,10).ClearContents()
$ws.Cells.Item(4,11).Value = "sense dades"
$ws.Cells.Item(4,12).Value = "sense dades sense dades"
$ws.Cells.Item(4,13).Value = "sense dades sense dades"
$ws.Cells.Item(4,14).Value = "sense dades sense dades"
$ws.Cells.Item(4,15).Value = "sense dades"

# Row 5
$ws.Cells.Item(5,4).Value = "'2026-02-07"
$ws.Cells.Item(5,5).Value = "2026-02-07 00:18:05"
$ws.Cells.Item(5,6).Value = "https://www.meteo.cat/observacions/xema/dades?codi=DJ&dia=2026-02-07T09:00Z"
$ws.Cells.Item(5,7).ClearContents()
$ws.Cells.Item(5,8).Value = "sense dades"
$ws.Cells.Item(5,9).Value = "sense dades"
$ws.Cells.Item(5,10).ClearContents()
$ws.Cells.Item(5,11).Value = "sense dades"
$ws.Cells.Item(5,12).Value = "sense dades sense dades"
$ws.Cells.Item(5,13).Value = "sense dades sense dades"
$ws.Cells.Item(5,14).Value = "sense dades sense dades"
$ws.Cells.Item(5,15).Value = "sense dades"

# Row 6
$ws.Cells.Item(6,4).Value = "'2026-02-07"
$ws.Cells.Item(6,5).Value = "2026-02-07 00:18:07"
$ws.Cells.Item(6,6).Value = "https://www.meteo.cat/observacions/xema/dades?codi=X4&dia=2026-02-07T09:00Z"
$ws.Cells.Item(6,7).ClearContents()
$ws.Cells.Item(6,8).Value = "sense dades"
$ws.Cells.Item(6,9).Value = "sense dades"
$ws.Cells.Item(6,10).ClearContents()
$ws.Cells.Item(6,11).Value = "sense dades"
$ws.Cells.Item(6,12).Value = "sense dades sense dades"
$ws.Cells.Item(6,13).Value = "sense dades sense dades"
$ws.Cells.Item(6,14).Value = "sense dades sense dades"
$ws.Cells.Item(6,15).Value = "sense dades"

# Row 7
$ws.Cells.Item(7,4).Value = "'2026-02-07"
$ws.Cells.Item(7,5).Value = "2026-02-07 00:18:09"
$ws.Cells.Item(7,6).Value = "https://www.meteo.cat/observacions/xema/dades?codi=D5&dia=2026-02-07T09:00Z"
$ws.Cells.Item(7,7).ClearContents()
$ws.Cells.Item(7,8).Value = "sense dades"
$ws.Cells.Item(7,9).Value = "sense dades"
$ws.Cells.Item(7,10).ClearContents()
$ws.Cells.Item(7,11).Value = "sense dades"
$ws.Cells.Item(7,12).Value = "sense dades sense dades"
$ws.Cells.Item(7,13).Value = "sense dades sense dades"
$ws.Cells.Item(7,14).Value = "sense dades sense dades"
$ws.Cells.Item(7,15).Value = "sense dades"

# Row 8
$ws.Cells.Item(8,4).Value = "'2026-02-07"
$ws.Cells.Item(8,5).Value = "2026-02-07 00:18:12"
$ws.Cells.Item(8,6).Value = "https://www.meteo.cat/observacions/xema/dades?codi=UN&dia=2026-02-07T09:00Z"
$ws.Cells.Item(8,7).ClearContents()
$ws.Cells.Item(8,8).Value = "sense dades"
$ws.Cells.Item(8,9).Value = "sense dades"
$ws.Cells.Item(8,10).ClearContents()
$ws.Cells.Item(8,11).Value = "sense dades"
$ws.Cells.Item(8,12).Value = "sense dades sense dades"
$ws.Cells.Item(8,13).Value = "sense dades sense dades"
$ws.Cells.Item(8,14).Value = "sense dades sense dades"
$ws.Cells.Item(8,15).Value = "sense dades"

# Row 9
$ws.Cells.Item(9,4).Value = "'2026-02-07"
$ws.Cells.Item(9,5).Value = "2026-02-07 00:18:14"
$ws.Cells.Item(9,6).Value = "https://www.meteo.cat/observacions/xema/dades?codi=MS&dia=2026-02-07T09:00Z"
$ws.Cells.Item(9,7).ClearContents()
$ws.Cells.Item(9,8).Value = "sense dades"
$ws.Cells.Item(9,9).Value = "sense dades"
$ws.Cells.Item(9,10).ClearContents()
$ws.Cells.Item(9,11).ClearContents()
$ws.Cells.Item(9,12).ClearContents()
$ws.Cells.Item(9,13).Value = "sense dades sense dades"
$ws.Cells.Item(9,14).Value = "sense dades sense dades"
$ws.Cells.Item(9,15).Value = "sense dades"

# Row 10
$ws.Cells.Item(10,4).Value = "'2026-02-07"
$ws.Cells.Item(10,5).Value = "2026-02-07 00:18:16"
$ws.Cells.Item(10,6).Value = "https://www.meteo.cat/observacions/xema/dades?codi=W1&dia=2026-02-07T09:00Z"
$ws.Cells.Item(10,7).ClearContents()
$ws.Cells.Item(10,8).Value = "sense dades"
$ws.Cells.Item(10,9).Value = "sense dades"
$ws.Cells.Item(10,10).ClearContents()
$ws.Cells.Item(10,11).ClearContents()
$ws.Cells.Item(10,12).ClearContents()
$ws.Cells.Item(10,13).Value = "sense dades sense dades"
$ws.Cells.Item(10,14).Value = "sense dades sense dades"
$ws.Cells.Item(10,15).Value = "sense dades"

# Row 11
$ws.Cells.Item(11,4).Value = "'2026-02-07"
$ws.Cells.Item(11,5).Value = "2026-02-07 00:18:18"
$ws.Cells.Item(11,6).Value = "https://www.meteo.cat/observacions/xema/dades?codi=DP&dia=2026-02-07T09:00Z"
$ws.Cells.Item(11,7).Value = "sense dades"
$ws.Cells.Item(11,8).Value = "sense dades"
$ws.Cells.Item(11,9).Value = "sense dades"
$ws.Cells.Item(11,10).ClearContents()
$ws.Cells.Item(11,11).Value = "sense dades"
$ws.Cells.Item(11,12).Value = "sense dades sense dades"
$ws.Cells.Item(11,13).Value = "sense dades sense dades"
$ws.Cells.Item(11,14).Value = "sense dades sense dades"
$ws.Cells.Item(11,15).Value = "sense dades"

# Row 12
$ws.Cells.Item(12,4).Value = "'2026-02-07"
$ws.Cells.Item(12,5).Value = "2026-02-07 00:18:21"
$ws.Cells.Item(12,6).Value = "https://www.meteo.cat/observacions/xema/dades?codi=XL&dia=2026-02-07T09:00Z"
$ws.Cells.Item(12,7).ClearContents()
$ws.Cells.Item(12,8).Value = "sense dades"
$ws.Cells.Item(12,9).Value = "sense dades"
$ws.Cells.Item(12,10).ClearContents()
$ws.Cells.Item(12,11).Value = "sense dades"
$ws.Cells.Item(12,12).Value = "sense dades sense dades"
$ws.Cells.Item(12,13).Value = "sense dades sense dades"
$ws.Cells.Item(12,14).Value = "sense dades sense dades"
$ws.Cells.Item(12,15).Value = "sense dades"

# Row 13
$ws.Cells.Item(13,4).Value = "'2026-02-07"
$ws.Cells.Item(13,5).Value = "2026-02-07 00:18:23"
$ws.Cells.Item(13,6).Value = "https://www.meteo.cat/observacions/xema/dades?codi=VZ&dia=2026-02-07T09:00Z"
$ws.Cells.Item(13,7).ClearContents()
$ws.Cells.Item(13,8).Value = "sense dades"
$ws.Cells.Item(13,9).Value = "sense dades"
$ws.Cells.Item(13,10).ClearContents()
$ws.Cells.Item(13,11).ClearContents()
$ws.Cells.Item(13,12).ClearContents()
$ws.Cells.Item(13,13).Value = "sense dades sense dades"
$ws.Cells.Item(13,14).Value = "sense dades sense dades"
$ws.Cells.Item(13,15).Value = "sense dades"

# Row 14
$ws.Cells.Item(14,4).Value = "'2026-02-07"
$ws.Cells.Item(14,5).Value = "2026-02-07 00:18:25"
$ws.Cells.Item(14,6).Value = "https://www.meteo.cat/observacions/xema/dades?codi=Z7&dia=2026-02-07T09:00Z"
$ws.Cells.Item(14,7).Value = "sense dades"
$ws.Cells.Item(14,8).Value = "sense dades"
$ws.Cells.Item(14,9).Value = "sense dades"
$ws.Cells.Item(14,10).ClearContents()
$ws.Cells.Item(14,11).Value = "sense dades"
$ws.Cells.Item(14,12).Value = "sense dades sense dades"
$ws.Cells.Item(14,13).Value = "sense dades sense dades"
$ws.Cells.Item(14,14).Value = "sense dades sense dades"
$ws.Cells.Item(14,15).Value = "sense dades"

# Row 15
$ws.Cells.Item(15,4).Value = "'2026-02-07"
$ws.Cells.Item(15,5).Value = "2026-02-07 00:18:28"
$ws.Cells.Item(15,6).Value = "https://www.meteo.cat/observacions/xema/dades?codi=XJ&dia=2026-02-07T09:00Z"
$ws.Cells.Item(15,7).ClearContents()
$ws.Cells.Item(15,8).Value = "sense dades"
$ws.Cells.Item(15,9).Value = "sense dades"
$ws.Cells.Item(15,10).ClearContents()
$ws.Cells.Item(15,11).Value = "sense dades"
$ws.Cells.Item(15,12).Value = "sense dades sense dades"
$ws.Cells.Item(15,13).Value = "sense dades sense dades"
$ws.Cells.Item(15,14).Value = "sense dades sense dades"
$ws.Cells.Item(15,15).Value = "sense dades"

# Row 16
$ws.Cells.Item(16,4).Value = "'2026-02-07"
$ws.Cells.Item(16,5).Value = "2026-02-07 00:18:30"
$ws.Cells.Item(16,6).Value = "https://www.meteo.cat/observacions/xema/dades?codi=YU&dia=2026-02-07T09:00Z"
$ws.Cells.Item(16,7).ClearContents()
$ws.Cells.Item(16,8).Value = "sense dades"
$ws.Cells.Item(16,9).Value = "sense dades"
$ws.Cells.Item(16,10).ClearContents()
$ws.Cells.Item(16,11).Value = "sense dades"
$ws.Cells.Item(16,12).Value = "sense dades sense dades"
$ws.Cells.Item(16,13).Value = "sense dades sense dades"
$ws.Cells.Item(16,14).Value = "sense dades sense dades"
$ws.Cells.Item(16,15).Value = "sense dades"

# Row 17
$ws.Cells.Item(17,4).Value = "'2026-02-07"
$ws.Cells.Item(17,5).Value = "2026-02-07 00:18:33"
$ws.Cells.Item(17,6).Value = "https://www.meteo.cat/observacions/xema/dades?codi=CD&dia=2026-02-07T09:00Z"
$ws.Cells.Item(17,7).ClearContents()
$ws.Cells.Item(17,8).Value = "sense dades"
$ws.Cells.Item(17,9).Value = "sense dades"
$ws.Cells.Item(17,10).ClearContents()
$ws.Cells.Item(17,11).Value = "sense dades"
$ws.Cells.Item(17,12).Value = "sense dades sense dades"
$ws.Cells.Item(17,13).Value = "sense dades sense dades"
$ws.Cells.Item(17,14).Value = "sense dades sense dades"
$ws.Cells.Item(17,15).Value = "sense dades"

# Row 18
$ws.Cells.Item(18,4).Value = "'2026-02-07"
$ws.Cells.Item(18,5).Value = "2026-02-07 00:18:35"
$ws.Cells.Item(18,6).Value = "https://www.meteo.cat/observacions/xema/dades?codi=Z2&dia=2026-02-07T09:00Z"
$ws.Cells.Item(18,7).Value = "sense dades"
$ws.Cells.Item(18,8).Value = "sense dades"
$ws.Cells.Item(18,9).Value = "sense dades"
$ws.Cells.Item(18,10).ClearContents()
$ws.Cells.Item(18,11).Value = "sense dades"
$ws.Cells.Item(18,12).Value = "sense dades sense dades"
$ws.Cells.Item(18,13).Value = "sense dades sense dades"
$ws.Cells.Item(18,14).Value = "sense dades sense dades"
$ws.Cells.Item(18,15).Value = "sense dades"

# Row 19
$ws.Cells.Item(19,4).Value = "'2026-02-07"
$ws.Cells.Item(19,5).Value = "2026-02-07 00:18:38"
$ws.Cells.Item(19,6).Value = "https://www.meteo.cat/observacions/xema/dades?codi=VK&dia=2026-02-07T09:00Z"
$ws.Cells.Item(19,7).ClearContents()
$ws.Cells.Item(19,8).Value = "sense dades"
$ws.Cells.Item(19,9).Value = "sense dades"
$ws.Cells.Item(19,10).ClearContents()
$ws.Cells.Item(19,11).Value = "sense dades"
$ws.Cells.Item(19,12).Value = "sense dades sense dades"
$ws.Cells.Item(19,13).Value = "sense dades sense dades"
$ws.Cells.Item(19,14).Value = "sense dades sense dades"
$ws.Cells.Item(19,15).Value = "sense dades"

# Row 20
$ws.Cells.Item(20,4).Value = "'2026-02-07"
$ws.Cells.Item(20,5).Value = "2026-02-07 00:18:40"
$ws.Cells.Item(20,6).Value = "https://www.meteo.cat/observacions/xema/dades?codi=Z3&dia=2026-02-07T09:00Z"
$ws.Cells.Item(20,7).Value = "sense dades"
$ws.Cells.Item(20,8).Value = "sense dades"
$ws.Cells.Item(20,9).Value = "sense dades"
$ws.Cells.Item(20,10).ClearContents()
$ws.Cells.Item(20,11).Value = "sense dades"
$ws.Cells.Item(20,12).Value = "sense dades sense dades"
$ws.Cells.Item(20,13).Value = "sense dades sense dades"
$ws.Cells.Item(20,14).Value = "sense dades sense dades"
$ws.Cells.Item(20,15).Value = "sense dades"

# Row 21
$ws.Cells.Item(21,4).Value = "'2026-02-07"
$ws.Cells.Item(21,5).Value = "2026-02-07 00:18:42"
$ws.Cells.Item(21,6).Value = "https://www.meteo.cat/observacions/xema/dades?codi=YB&dia=2026-02-07T09:00Z"
$ws.Cells.Item(21,7).ClearContents()
$ws.Cells.Item(21,8).Value = "sense dades"
$ws.Cells.Item(21,9).Value = "sense dades"
$ws.Cells.Item(21,10).ClearContents()
$ws.Cells.Item(21,11).Value = "sense dades"
$ws.Cells.Item(21,12).Value = "sense dades sense dades"
$ws.Cells.Item(21,13).Value = "sense dades sense dades"
$ws.Cells.Item(21,14).Value = "sense dades sense dades"
$ws.Cells.Item(21,15).Value = "sense dades"

# Row 22
$ws.Cells.Item(22,4).Value = "'2026-02-07"
$ws.Cells.Item(22,5).Value = "2026-02-07 00:18:44"
$ws.Cells.Item(22,6).Value = "https://www.meteo.cat/observacions/xema/dades?codi=YP&dia=2026-02-07T09:00Z"
$ws.Cells.Item(22,7).ClearContents()
$ws.Cells.Item(22,8).Value = "sense dades"
$ws.Cells.Item(22,9).Value = "sense dades"
$ws.Cells.Item(22,10).ClearContents()
$ws.Cells.Item(22,11).Value = "sense dades"
$ws.Cells.Item(22,12).Value = "sense dades sense dades"
$ws.Cells.Item(22,13).Value = "sense dades sense dades"
$ws.Cells.Item(22,14).Value = "sense dades sense dades"
$ws.Cells.Item(22,15).Value = "sense dades"

# Row 23
$ws.Cells.Item(23,4).Value = "'2026-02-07"
$ws.Cells.Item(23,5).Value = "2026-02-07 00:18:47"
$ws.Cells.Item(23,6).Value = "https://www.meteo.cat/observacions/xema/dades?codi=J5&dia=2026-02-07T09:00Z"
$ws.Cells.Item(23,7).ClearContents()
$ws.Cells.Item(23,8).Value = "sense dades"
$ws.Cells.Item(23,9).Value = "sense dades"
$ws.Cells.Item(23,10).ClearContents()
$ws.Cells.Item(23,11).Value = "sense dades"
$ws.Cells.Item(23,12).Value = "sense dades sense dades"
$ws.Cells.Item(23,13).Value = "sense dades sense dades"
$ws.Cells.Item(23,14).Value = "sense dades sense dades"
$ws.Cells.Item(23,15).Value = "sense dades"

# Row 24
$ws.Cells.Item(24,4).Value = "'2026-02-07"
$ws.Cells.Item(24,5).Value = "2026-02-07 00:18:49"
$ws.Cells.Item(24,6).Value = "https://www.meteo.cat/observacions/xema/dades?codi=D6&dia=2026-02-07T09:00Z"
$ws.Cells.Item(24,7).ClearContents()
$ws.Cells.Item(24,8).Value = "sense dades"
$ws.Cells.Item(24,9).Value = "sense dades"
$ws.Cells.Item(24,10).ClearContents()
$ws.Cells.Item(24,11).Value = "sense dades"
$ws.Cells.Item(24,12).Value = "sense dades sense dades"
$ws.Cells.Item(24,13).Value = "sense dades sense dades"
$ws.Cells.Item(24,14).Value = "sense dades sense dades"
$ws.Cells.Item(24,15).Value = "sense dades"

# Row 25
$ws.Cells.Item(25,4).Value = "'2026-02-07"
$ws.Cells.Item(25,5).Value = "2026-02-07 00:18:51"
$ws.Cells.Item(25,6).Value = "https://www.meteo.cat/observacions/xema/dades?codi=YA&dia=2026-02-07T09:00Z"
$ws.Cells.Item(25,7).ClearContents()
$ws.Cells.Item(25,8).Value = "sense dades"
$ws.Cells.Item(25,9).Value = "sense dades"
$ws.Cells.Item(25,10).ClearContents()
$ws.Cells.Item(25,11).Value = "sense dades"
$ws.Cells.Item(25,12).Value = "sense dades sense dades"
$ws.Cells.Item(25,13).Value = "sense dades sense dades"
$ws.Cells.Item(25,14).Value = "sense dades sense dades"
$ws.Cells.Item(25,15).Value = "sense dades"

# Row 26
$ws.Cells.Item(26,4).Value = "'2026-02-07"
$ws.Cells.Item(26,5).Value = "2026-02-07 00:18:54"
$ws.Cells.Item(26,6).Value = "https://www.meteo.cat/observacions/xema/dades?codi=DG&dia=2026-02-07T09:00Z"
$ws.Cells.Item(26,7).Value = "sense dades"
$ws.Cells.Item(26,8).Value = "sense dades"
$ws.Cells.Item(26,9).Value = "sense dades"
$ws.Cells.Item(26,10).ClearContents()
$ws.Cells.Item(26,11).Value = "sense dades"
$ws.Cells.Item(26,12).Value = "sense dades sense dades"
$ws.Cells.Item(26,13).Value = "sense dades sense dades"
$ws.Cells.Item(26,14).Value = "sense dades sense dades"
$ws.Cells.Item(26,15).Value = "sense dades"

# Row 27
$ws.Cells.Item(27,4).Value = "'2026-02-07"
$ws.Cells.Item(27,5).Value = "2026-02-07 00:18:56"
$ws.Cells.Item(27,6).Value = "https://www.meteo.cat/observacions/xema/dades?codi=D4&dia=2026-02-07T09:00Z"
$ws.Cells.Item(27,7).ClearContents()
$ws.Cells.Item(27,8).Value = "sense dades"
$ws.Cells.Item(27,9).Value = "sense dades"
$ws.Cells.Item(27,10).ClearContents()
$ws.Cells.Item(27,11).Value = "sense dades"
$ws.Cells.Item(27,12).Value = "sense dades sense dades"
$ws.Cells.Item(27,13).Value = "sense dades sense dades"
$ws.Cells.Item(27,14).Value = "sense dades sense dades"
$ws.Cells.Item(27,15).Value = "sense dades"

# Row 28
$ws.Cells.Item(28,4).Value = "'2026-02-07"
$ws.Cells.Item(28,5).Value = "2026-02-07 00:18:59"
$ws.Cells.Item(28,6).Value = "https://www.meteo.cat/observacions/xema/dades?codi=CI&dia=2026-02-07T09:00Z"
$ws.Cells.Item(28,7).ClearContents()
$ws.Cells.Item(28,8).Value = "sense dades"
$ws.Cells.Item(28,9).Value = "sense dades"
$ws.Cells.Item(28,10).ClearContents()
$ws.Cells.Item(28,11).ClearContents()
$ws.Cells.Item(28,12).Value = "sense dades sense dades"
$ws.Cells.Item(28,13).Value = "sense dades sense dades"
$ws.Cells.Item(28,14).Value = "sense dades sense dades"
$ws.Cells.Item(28,15).Value = "sense dades"

# Row 29
$ws.Cells.Item(29,4).Value = "'2026-02-07"
$ws.Cells.Item(29,5).Value = "2026-02-07 00:19:01"
$ws.Cells.Item(29,6).Value = "https://www.meteo.cat/observacions/xema/dades?codi=XS&dia=2026-02-07T09:00Z"
$ws.Cells.Item(29,7).ClearContents()
$ws.Cells.Item(29,8).Value = "sense dades"
$ws.Cells.Item(29,9).Value = "sense dades"
$ws.Cells.Item(29,10).ClearContents()
$ws.Cells.Item(29,11).Value = "sense dades"
$ws.Cells.Item(29,12).Value = "sense dades sense dades"
$ws.Cells.Item(29,13).Value = "sense dades sense dades"
$ws.Cells.Item(29,14).Value = "sense dades sense dades"
$ws.Cells.Item(29,15).Value = "sense dades"

# Row 30
$ws.Cells.Item(30,4).Value = "'2026-02-07"
$ws.Cells.Item(30,5).Value = "2026-02-07 00:19:03"
$ws.Cells.Item(30,6).Value = "https://www.meteo.cat/observacions/xema/dades?codi=ZC&dia=2026-02-07T09:00Z"
$ws.Cells.Item(30,7).Value = "sense dades"
$ws.Cells.Item(30,8).Value = "sense dades"
$ws.Cells.Item(30,9).Value = "sense dades"
$ws.Cells.Item(30,10).ClearContents()
$ws.Cells.Item(30,11).Value = "sense dades"
$ws.Cells.Item(30,12).Value = "sense dades sense dades"
$ws.Cells.Item(30,13).Value = "sense dades sense dades"
$ws.Cells.Item(30,14).Value = "sense dades sense dades"
$ws.Cells.Item(30,15).Value = "sense dades"

# Row 31
$ws.Cells.Item(31,4).Value = "'2026-02-07"
$ws.Cells.Item(31,5).Value = "2026-02-07 00:19:06"
$ws.Cells.Item(31,6).Value = "https://www.meteo.cat/observacions/xema/dades?codi=XH&dia=2026-02-07T09:00Z"
$ws.Cells.Item(31,7).Value = "sense dades"
$ws.Cells.Item(31,8).Value = "sense dades"
$ws.Cells.Item(31,9).Value = "sense dades"
$ws.Cells.Item(31,10).ClearContents()
$ws.Cells.Item(31,11).ClearContents()
$ws.Cells.Item(31,12).ClearContents()
$ws.Cells.Item(31,13).Value = "sense dades sense dades"
$ws.Cells.Item(31,14).Value = "sense dades sense dades"
$ws.Cells.Item(31,15).Value = "sense dades"

# Row 32
$ws.Cells.Item(32,4).Value = "'2026-02-07"
$ws.Cells.Item(32,5).Value = "2026-02-07 00:19:08"
$ws.Cells.Item(32,6).Value = "https://www.meteo.cat/observacions/xema/dades?codi=XE&dia=2026-02-07T09:00Z"
$ws.Cells.Item(32,7).ClearContents()
$ws.Cells.Item(32,8).Value = "sense dades"
$ws.Cells.Item(32,9).Value = "sense dades"
$ws.Cells.Item(32,10).ClearContents()
$ws.Cells.Item(32,11).Value = "sense dades"
$ws.Cells.Item(32,12).Value = "sense dades sense dades"
$ws.Cells.Item(32,13).Value = "sense dades sense dades"
$ws.Cells.Item(32,14).Value = "sense dades sense dades"
$ws.Cells.Item(32,15).Value = "sense dades"

# Row 33
$ws.Cells.Item(33,4).Value = "'2026-02-07"
$ws.Cells.Item(33,5).Value = "2026-02-07 00:19:10"
$ws.Cells.Item(33,6).Value = "https://www.meteo.cat/observacions/xema/dades?codi=UE&dia=2026-02-07T09:00Z"
$ws.Cells.Item(33,7).ClearContents()
$ws.Cells.Item(33,8).Value = "sense dades"
$ws.Cells.Item(33,9).Value = "sense dades"
$ws.Cells.Item(33,10).ClearContents()
$ws.Cells.Item(33,11).ClearContents()
$ws.Cells.Item(33,12).ClearContents()
$ws.Cells.Item(33,13).Value = "sense dades sense dades"
$ws.Cells.Item(33,14).Value = "sense dades sense dades"
$ws.Cells.Item(33,15).Value = "sense dades"

# Row 34
$ws.Cells.Item(34,4).Value = "'2026-02-07"
$ws.Cells.Item(34,5).Value = "2026-02-07 00:19:13"
$ws.Cells.Item(34,6).Value = "https://www.meteo.cat/observacions/xema/dades?codi=XO&dia=2026-02-07T09:00Z"
$ws.Cells.Item(34,7).ClearContents()
$ws.Cells.Item(34,8).Value = "sense dades"
$ws.Cells.Item(34,9).Value = "sense dades"
$ws.Cells.Item(34,10).ClearContents()
$ws.Cells.Item(34,11).Value = "sense dades"
$ws.Cells.Item(34,12).Value = "sense dades sense dades"
$ws.Cells.Item(34,13).Value = "sense dades sense dades"
$ws.Cells.Item(34,14).Value = "sense dades sense dades"
$ws.Cells.Item(34,15).Value = "sense dades"

# Row 35
$ws.Cells.Item(35,4).Value = "'2026-02-07"
$ws.Cells.Item(35,5).Value = "2026-02-07 00:19:15"
$ws.Cells.Item(35,6).Value = "https://www.meteo.cat/observacions/xema/dades?codi=VS&dia=2026-02-07T09:00Z"
$ws.Cells.Item(35,7).Value = "sense dades"
$ws.Cells.Item(35,8).Value = "sense dades"
$ws.Cells.Item(35,9).Value = "sense dades"
$ws.Cells.Item(35,10).ClearContents()
$ws.Cells.Item(35,11).Value = "sense dades"
$ws.Cells.Item(35,12).Value = "sense dades sense dades"
$ws.Cells.Item(35,13).Value = "sense dades sense dades"
$ws.Cells.Item(35,14).Value = "sense dades sense dades"
$ws.Cells.Item(35,15).Value = "sense dades"

# Row 36
$ws.Cells.Item(36,4).Value = "'2026-02-07"
$ws.Cells.Item(36,5).Value = "2026-02-07 00:19:17"
$ws.Cells.Item(36,6).Value = "https://www.meteo.cat/observacions/xema/dades?codi=D7&dia=2026-02-07T09:00Z"
$ws.Cells.Item(36,7).ClearContents()
$ws.Cells.Item(36,8).Value = "sense dades"
$ws.Cells.Item(36,9).Value = "sense dades"
$ws.Cells.Item(36,10).ClearContents()
$ws.Cells.Item(36,11).Value = "sense dades"
$ws.Cells.Item(36,12).Value = "sense dades sense dades"
$ws.Cells.Item(36,13).Value = "sense dades sense dades"
$ws.Cells.Item(36,14).Value = "sense dades sense dades"
$ws.Cells.Item(36,15).Value = "sense dades"
